# Gantt chart update (4/22/2020 mentor discussion).
#
# The original "200413" sheet is duplicated to a new "200422" tab that keeps
# (most of) the old data with a handful of edits/new rows, while the
# original "200413" sheet is rewritten in place with a new "Category"
# column and a few renamed tasks.

$wb = $excel.ActiveWorkbook
$orig = $wb.Worksheets.Item(1)
$orig.Name = "200413"

# --- 1. Create the new "200422" tab right after "200413", as a copy of ---
# --- the ORIGINAL (pre-edit) 200413 data, then apply its own edits.    ---
$ws2 = $wb.Worksheets.Add($null, $orig)
$ws2.Name = "200422"

$ws2.Columns.Item(1).ColumnWidth = 43.7265625
$ws2.Columns.Item(2).ColumnWidth = 9.31640625
$ws2.Columns.Item(3).ColumnWidth = 8.90625
$ws2.Columns.Item(4).ColumnWidth = 9.2265625

$rows2 = @(
    @{ Row=1;  A="Task";                                                         B="Category";     C="Start Date"; D="End Date" },
    @{ Row=2;  A="Descriptor/Property correlation";                              B="Initial data"; C=43934;        D=43943 },
    @{ Row=3;  A="ML exploration (NN, RFR, KRR)";                                B="Initial data"; C=43936;        D=43943 },
    @{ Row=4;  A="Train models with DFT data";                                   B="Initial data"; C=43936;        D=43943 },
    @{ Row=5;  A="Clean DFT data: outliers, normalize";                          B="Total data";   C=43943;        D=43950; Wrap=$true },
    @{ Row=6;  A="Descriptor importance: remove unecessary";                     B="Total data";   C=43945;        D=43950; Wrap=$true },
    @{ Row=7;  A="Train models with new DFT data (NN, RFR, LASSO, GPR)";         B="Total data";   C=43945;        D=43957; Wrap=$true; Tall=$true;  BuiltinDate=$true },
    @{ Row=8;  A="Test models with mixed alloy data";                           B="Total data";   C=43964;        D=43971; BuiltinDate=$true },
    @{ Row=9;  A="Brainstorm more descriptors";                                 B="Lit Review";   C=43966;        D=43973; BDateFmt=$true; BuiltinDate=$true },
    @{ Row=10; A="Expected improvement:informed selection of new training data"; B="Bonus";        C=43971;        D=43979; Wrap=$true; Tall=$true; BDateFmt=$true; BuiltinDate=$true },
    @{ Row=11; A="Have model completed";                                        B="Overall";      C=43987;        D=44001 }
)

foreach ($r in $rows2) {
    $row = $r.Row
    $ws2.Cells.Item($row, 1).Value = $r.A
    $ws2.Cells.Item($row, 2).Value = $r.B

    if ($row -eq 1) {
        $ws2.Cells.Item($row, 3).Value = $r.C
        $ws2.Cells.Item($row, 4).Value = $r.D
    } else {
        $ws2.Cells.Item($row, 3).Value = $r.C
        $ws2.Cells.Item($row, 3).NumberFormat = "mm/dd/yy;@"
        $ws2.Cells.Item($row, 4).Value = $r.D

        if (-not $r.BuiltinDate) {
            $ws2.Cells.Item($row, 4).NumberFormat = "mm/dd/yy;@"
        }
    }

    if ($r.Wrap) {
        $ws2.Cells.Item($row, 1).WrapText = $true
    }
    if ($r.BDateFmt) {
        $ws2.Cells.Item($row, 2).NumberFormat = "mm/dd/yy;@"
    }
    if ($r.Tall) {
        $ws2.Rows.Item($row).RowHeight = 29.5
    }
}

# D7, D8, D9 & D10 use the plain builtin date style (same one already used
# by the pre-existing "have model completed" row further down / row 6 of
# the original sheet) instead of the custom "mm/dd/yy;@" style.
$orig.Range("C6").Copy() | Out-Null
$ws2.Range("D7:D10").PasteSpecial(-4122) | Out-Null

$ws2.Range("A10").Select() | Out-Null

# --- 2. Rewrite the original "200413" sheet in place. ---------------------
$orig.Columns.Item(2).Insert() | Out-Null

$orig.Columns.Item(1).ColumnWidth = 43.7265625
$orig.Columns.Item(2).ColumnWidth = 9.31640625
$orig.Columns.Item(3).ColumnWidth = 8.90625
$orig.Columns.Item(4).ColumnWidth = 9.2265625
$orig.Columns.Item(5).ColumnWidth = 8.2265625

$orig.Cells.Item(1, 1).Value = "Task"
$orig.Cells.Item(1, 2).Value = "Category"
$orig.Cells.Item(1, 3).Value = "Start Date"
$orig.Cells.Item(1, 4).Value = "End Date"

$rows1 = @(
    @{ Row=2; A="Descriptor/Property correlation";  B="Initial data"; C=43934; D=43943 },
    @{ Row=3; A="ML exploration (NN, RF, KRR)";      B="Initial data"; C=43936; D=43950 },
    @{ Row=4; A="Train models with data";            B="Initial data"; C=43936; D=43950 },
    @{ Row=5; A="Have model completed";              B="Overall";      C=43987; D=44001 }
)

foreach ($r in $rows1) {
    $row = $r.Row
    $orig.Cells.Item($row, 1).Value = $r.A
    $orig.Cells.Item($row, 2).Value = $r.B
    $orig.Cells.Item($row, 3).Value = $r.C
    $orig.Cells.Item($row, 3).NumberFormat = "mm/dd/yy;@"
    $orig.Cells.Item($row, 4).Value = $r.D
    $orig.Cells.Item($row, 4).NumberFormat = "mm/dd/yy;@"
    $orig.Cells.Item($row, 5).NumberFormat = "mm/dd/yy;@"
}

$orig.Range("A4").Select() | Out-Null

# --- 3. Workbook-level view state -----------------------------------------
$ws2.Activate()
